$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C140").Value = 45205
